$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> @(new Price/D value or $null, new Volume(1h)/E value or $null)
$updates = @{
    2 = @('42.972.23', '  +0.47%  ')
    3 = @('2.283.72', '  +0.42%  ')
    4 = @($null, '  -0.11%  ')
    5 = @('250.70', $null)
    6 = @('0.634', '  +1.10%  ')
    7 = @('78.11', '  +8.55%  ')
    8 = @($null, '  -0.01%  ')
    9 = @('0.657', '  +1.28%  ')
    10 = @('41.26', '  +6.81%  ')
    11 = @('0.0976', '  +1.16%  ')
    12 = @('7.37', '  +0.05%  ')
    13 = @('0.105', '  -0.31%  ')
    14 = @('2.624.30', '  +0.33%  ')
    15 = @('15.06', '  +0.60%  ')
    16 = @('0.868', '  -1.66%  ')
    17 = @('2.283.48', '  +0.62%  ')
    18 = @('42.906.21', '  +0.48%  ')
    19 = @('0.0₃0996', '  -2.00%  ')
    20 = @($null, '  -1.12%  ')
    21 = @('72.26', '  -1.54%  ')
    22 = @('234.57', '  +0.18%  ')
    23 = @('2.17', '  +1.34%  ')
    24 = @($null, '  -3.45%  ')
    25 = @($null, '  +0.00%  ')
    26 = @('11.40', '  -1.72%  ')
    27 = @('2.35', '  -3.79%  ')
    28 = @($null, '  +2.17%  ')
    29 = @('167.89', '  -0.17%  ')
    30 = @('20.97', '  -0.17%  ')
    31 = @($null, '  -0.03%  ')
    32 = @('0.0858', '  +7.23%  ')
    33 = @('0.124', '  -4.37%  ')
    34 = @('30.51', '  -1.42%  ')
    35 = @($null, '  +0.41%  ')
    36 = @($null, '  -0.98%  ')
    37 = @('4.76', '  +0.04%  ')
    38 = @('0.0306', '  -1.85%  ')
    39 = @('13.89', '  +3.57%  ')
    40 = @($null, '  -2.15%  ')
    41 = @('5.89', '  +0.98%  ')
    42 = @('112.55', '  +18.12%  ')
    43 = @($null, '  -1.30%  ')
    44 = @('61.41', '  -0.50%  ')
    45 = @('8.92', '  -2.75%  ')
    46 = @('0.102', '  -1.18%  ')
    47 = @($null, '  +0.03%  ')
    48 = @('4.57', '  -7.99%  ')
    49 = @($null, '  -2.29%  ')
    50 = @($null, '  -1.85%  ')
    51 = @('4.26', '  +0.45%  ')
}

foreach ($row in ($updates.Keys | Sort-Object)) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $dCell = $ws.Range("D$row")
        # Force text storage so numeric-looking strings (e.g. "250.70") keep their
        # exact textual representation instead of being normalized as a number.
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.Style = "Normal"
    }
    if ($null -ne $eVal) {
        $ws.Range("E$row").Value = $eVal
    }
}
